# Daily attendance processing - 2025-12-09 17:02:26
# Normalizes the "Recorded By" (column G) values: when a cell lists two or
# more recorders separated by ", ", the first two entries are swapped
# (this reorders e.g. "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System"). The sentinel value
# "System, backup@backdoor.com" is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$firstRow = $usedRange.Row

# Locate the "Recorded By" column by scanning the header row instead of
# hardcoding a column index.
$headerRow = $firstRow
$colCount = $usedRange.Columns.Count
$firstCol = $usedRange.Column
$recordedByCol = 0
for ($c = 0; $c -lt $colCount; $c++) {
    $colIndex = $firstCol + $c
    $headerText = $ws.Cells.Item($headerRow, $colIndex).Text
    if ($headerText.Equals("Recorded By")) {
        $recordedByCol = $colIndex
        break
    }
}

if ($recordedByCol -eq 0) {
    $recordedByCol = 7  # fallback to column G
}

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstRow + $i
    if ($r -eq $headerRow) {
        continue
    }

    $cell = $ws.Cells.Item($r, $recordedByCol)
    $current = $cell.Text

    if ([string]::IsNullOrEmpty($current)) {
        continue
    }

    # Leave the known sentinel combination untouched.
    if ($current.Equals("System, backup@backdoor.com")) {
        continue
    }

    # Split on the ", " separator (substring, not character set) so that
    # tokens themselves are never broken apart.
    $commaIndex = $current.IndexOf(", ")
    if ($commaIndex -ge 0) {
        $firstPart = $current.Substring(0, $commaIndex)
        $rest = $current.Substring($commaIndex + 2)

        $secondCommaIndex = $rest.IndexOf(", ")
        if ($secondCommaIndex -ge 0) {
            $secondPart = $rest.Substring(0, $secondCommaIndex)
            $remainder = $rest.Substring($secondCommaIndex)
        } else {
            $secondPart = $rest
            $remainder = ""
        }

        $updated = $secondPart + ", " + $firstPart + $remainder
        if (-not $updated.Equals($current)) {
            $cell.Value2 = $updated
        }
    }
}
